# Fix Update after update or change student's info and fix SortbyID
#
# Adds two new "teacher" login rows (giaovien / giaovien1) to the
# UserData sheet, re-using the existing "123" password shared string
# (same value already used by the "nhi" / "nhi1" rows in B7:B8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New username cells (plain text -> new shared-string entries).
$ws.Range("A9").Value  = "giaovien"
$ws.Range("A10").Value = "giaovien1"

# New password cells: copy the existing B7 cell ("123", already a shared
# string) into B9/B10 instead of assigning the literal text "123" again,
# so Excel keeps it as text (shared string) rather than re-inferring it
# as a number, and so no new cell style gets allocated.
$ws.Range("B7").Copy($ws.Range("B9"))
$ws.Range("B7").Copy($ws.Range("B10"))
